$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap ASP.NET")

# Set the actual duration for week 2 (row 4, column E)
$ws.Range("E4").Value = 90

# Recalculate the workbook so the dependent formulas (F4:F12) and chart caches update
$excel.Calculate()

# Update the view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I5").Select()
